$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of "fullRNASEQ" -> "fullRNASeq" in the purpose column (E)
# for every data row (rows 2 through 24).
for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 5).Value = "fullRNASeq"
}
